$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, reusing the same formatting as the other
# header cells (e.g. G1: bold, centered, thin box border) by copying the
# existing header cell's format rather than building a new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 for the single data row.
$ws.Range("H2").Value = 0
